$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update matrices (F column) values for rows 2-13
$ws.Range("F2").Value = 14.35604799398173
$ws.Range("F3").Value = 13.17756464437572
$ws.Range("F4").Value = 8.158367614863963
$ws.Range("F5").Value = 7.132419507397405
$ws.Range("F6").Value = 6.306267974076017
$ws.Range("F7").Value = 6.305467982787811
$ws.Range("F8").Value = 5.443833869706829
$ws.Range("F9").Value = 5.210446373867417
$ws.Range("F10").Value = 5.049269166493271
$ws.Range("F11").Value = 3.295791998891051
$ws.Range("F12").Value = 1.043195326962711
$ws.Range("F13").Value = 0.2516358054655306

# Update index (B column) values for rows 8, 9, 10
$ws.Range("B8").Value = 30
$ws.Range("B9").Value = 32
$ws.Range("B10").Value = 33

# Update race (G column) values for rows 8, 9 (swap)
$ws.Range("G8").Value = "White"
$ws.Range("G9").Value = "Black or African American"
